$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: apply updated values per diff
$ws.Range("D2").Value = 44299
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("R2").Value = 'Provincia de Santiago'
$ws.Range("S2").Value = 2143

# Row 3: apply updated values per diff
$ws.Range("D3").Value = 44299
$ws.Range("M3").Value = 75
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("R3").Value = 'Provincia de Santiago'
$ws.Range("S3").Value = 1714

# Row 4: apply updated values per diff
$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("S4").Value = 1714

# Row 5: apply updated values per diff
$ws.Range("D5").Value = 44322
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("S5").Value = 1143

# Row 6: apply updated values per diff
$ws.Range("D6").Value = 44300
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 2143

# Row 7: apply updated values per diff
$ws.Range("D7").Value = 44300
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 1714

# Row 8: apply updated values per diff
$ws.Range("D8").Value = 44302
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 2143

# Row 9: apply updated values per diff
$ws.Range("D9").Value = 44302
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 10: apply updated values per diff
$ws.Range("D10").Value = 44980
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 2286

# Row 11: apply updated values per diff
$ws.Range("D11").Value = 44980
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 13000
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 1857

# Row 12: apply updated values per diff
$ws.Range("D12").Value = 44320
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 1714

# Row 13: apply updated values per diff
$ws.Range("D13").Value = 44320
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("S13").Value = 1143

# Row 14: apply updated values per diff
$ws.Range("D14").Value = 44971
$ws.Range("M14").Value = 25
$ws.Range("Q14").Value = '$/bandeja 5 kilos'
$ws.Range("S14").Value = 3000
$ws.Range("T14").Value = 5

# Row 15: apply updated values per diff
$ws.Range("D15").Value = 44292
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("S15").Value = 2286

# Row 16: apply updated values per diff
$ws.Range("D16").Value = 44292
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 30
$ws.Range("Q16").Value = '$/bandeja 7 kilos'
$ws.Range("S16").Value = 2143
$ws.Range("T16").Value = 7
